$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ALLEGRETTO-LTE (B7981028) days remaining: 13 -> 12
$ws.Range("B7").Value = 12

# REMASTER (CLOU) days remaining: 33 -> 32
$ws.Range("B9").Value = 32

# COLO-PREVENT progress: new value 100
$ws.Range("C10").Value = 100
